$wb = $excel.ActiveWorkbook

# The "survey" sheet holds the note text that references the refrigerator id.
$survey = $wb.Worksheets.Item("survey")

# Update the handlebars expression to expose the data model (data.refrigerator_id)
# instead of the bare refrigerator_id.
$survey.Range("D2").Value = "Refrigerator id: {{data.refrigerator_id}}"

# Make the edited cell/sheet the active selection, mirroring the manual edit.
$survey.Activate()
$survey.Range("D2").Select()
